$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph (the 2nd paragraph
#    of the document, right after the H1 title).
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Blazing Bull for Free –
#    Exciting High Volatility Slot" right before the very last
#    paragraph (the former "Prompt: ..." paragraph).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertionPoint = $lastPara.Range
$insertionPoint.Collapse(1)              # wdCollapseStart
$insertionPoint.InsertParagraphBefore()

$newCount = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($newCount - 1)
$newRange = $newPara.Range

$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Blazing Bull for Free &#8211; Exciting High Volatility Slot</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRange.InsertXML($titleXml)

# ------------------------------------------------------------------
# 3) Replace the text of the final paragraph (formerly the "Prompt:
#    ..." image-generation prompt) with the meta-description copy,
#    keeping its existing (italic) run formatting intact.
# ------------------------------------------------------------------
$finalCount = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($finalCount)
$finalRange = $finalPara.Range

$oldText = "Prompt: Create a cartoon-style feature image for Blazing Bull that features a happy Maya warrior with glasses. The image should capture the essence of the game by including elements like flames, animals, and the slot reel. The Maya warrior should be holding a golden bull coin, surrounded by animals like bears and mountain lions. Flames should be visible in the background to represent the game's title " + [char]34 + "Blazing Bull." + [char]34 + " The warrior should be positioned in a way that suggests they are about to spin the slot reel. The image should be colorful and eye-catching, with a mystical aura surrounding it, similar to the game's aesthetic."
$newText = "Experience high volatility gameplay with Blazing Bull. Trigger bonus modes for free spins and multipliers. Play for free today."

$finalRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
